# Requirements.xlsx update: add a "Data" sheet documenting sensor fields,
# drop the fire-detection requirements/spec/functional rows, tweak the
# 12-hour transmission wording, and add a GPS/cellular technical spec row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Technical Specs sheet: remove the two fire-detection rows, add a new
#    "Cellular Base station" row, and renumber column A.
# ---------------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("Technical Specs")
$wsTech.Rows.Item(5).Delete()
$wsTech.Rows.Item(5).Delete()
$wsTech.Rows.Item(9).Insert()
$wsTech.Range("A4").Value = 4
$wsTech.Range("A5").Value = 5
$wsTech.Range("A6").Value = 6
$wsTech.Range("A7").Value = 7
$wsTech.Range("A8").Value = 8
$wsTech.Range("A9").Value = 9
$wsTech.Range("B9").Value = "The device will use the Cellular Base station on determining position in case the GPS cannot get a fix"
$wsTech.Range("A10").Value = 10
$wsTech.Columns.Item(2).ColumnWidth = 93.7109375

# ---------------------------------------------------------------------
# 2) Functional Req sheet: drop the fire-alert row, tighten the 12-hour
#    transmission wording.
# ---------------------------------------------------------------------
$wsFunc = $wb.Worksheets.Item("Functional Req")
$wsFunc.Rows.Item(4).Delete()
$wsFunc.Range("A4").Value = 4
$wsFunc.Range("A5").Value = 5
$wsFunc.Range("A6").Value = 6
$wsFunc.Range("B1").Value = "When 12 hours is elapsed since the last transmission of the sensor data (trash level, orientation,location), the device will transmit the current sensor data."
$wsFunc.Rows.Item(1).RowHeight = 30

# ---------------------------------------------------------------------
# 3) Insert the new "Data" sheet right before "Non-Functional Req" and
#    fill in the sensor-field documentation table.
# ---------------------------------------------------------------------
$wsNonFunc = $wb.Worksheets.Item("Non-Functional Req")
$wsData = $wb.Worksheets.Add($wsNonFunc)
$wsData.Name = "Data"

$wsData.Range("A1").Value = "Name"
$wsData.Range("B1").Value = "Type"
$wsData.Range("C1").Value = "Range"
$wsData.Range("D1").Value = "Unit"

$wsData.Range("A2").Value = "Fill Level"
$wsData.Range("B2").Value = "unsigned int"
$wsData.Range("C2").Value = "0 - 100"
$wsData.Range("D2").Value = "%"

$wsData.Range("A3").Value = "Orientation"
$wsData.Range("B3").Value = "unsigned int"
$wsData.Range("C3").Value = "0 - 360"
$wsData.Range("D3").Value = "Degrees"

$wsData.Range("A4").Value = "Location_GPS"
$wsData.Range("B4").Value = "string"

$wsData.Range("A5").Value = "Location_GSM"
$wsData.Range("B5").Value = "string"

$wsData.Range("A6").Value = "Calibration Status"
$wsData.Range("B6").Value = "bool"
$wsData.Range("C6").Value = "true,false"

$wsData.Range("A7").Value = "Signal Level"
$wsData.Range("B7").Value = "int"
$wsData.Range("D7").Value = "dbm"

$wsData.Columns.Item(1).ColumnWidth = 16.28515625
$wsData.Columns.Item(2).ColumnWidth = 16.42578125
$wsData.Columns.Item(3).ColumnWidth = 13.85546875
$wsData.Columns.Item(4).ColumnWidth = 12.5703125

# ---------------------------------------------------------------------
# 4) Update selections on each sheet (matches the saved UI state), then
#    leave "Functional Req" as the active tab, as before.
# ---------------------------------------------------------------------
$wsTech.Activate()
$wsTech.Range("B9").Select()

$wsData.Activate()
$wsData.Range("B15").Select()

$wsNonFunc.Activate()
$wsNonFunc.Range("O20").Select()

$wsFunc.Activate()
$wsFunc.Range("G12").Select()
